# Apply the tweak described by the commit "tweaked app for heroku":
# A new blank column ("Unnamed: 0") is inserted before the old EMPLOYEE_ID
# column, shifting EMPLOYEE_ID..DEPARTMENT one column to the right
# (D:J -> E:K). The pivot was also refreshed against an updated employee
# lookup table, so MANAGER_ID (col F), PROCESS (col J) and, for a couple
# of rows, DEPARTMENT (col K) come out with new values rather than just
# being shifted copies of the old data - those are written explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts the existing D:J data to E:K
# and automatically grows the sheet dimension from A1:J15 to A1:K15.
$ws.Columns("D:D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Unnamed: 0"

# Refreshed MANAGER_ID (F), PROCESS (J) and DEPARTMENT (K) values per row.
$managerId  = @{ 2="O254"; 3="O50"; 4="O72"; 5="O239"; 6="O254"; 7="O50"; 8="O50"; 9="O72"; 10="O72"; 11="O72"; 12="O72"; 13="O72"; 14="O72"; 15="O50" }
$process    = @{ 2="MUTHOOT"; 3="IDFC"; 4="L&T"; 5="IDFC"; 6="IDFC"; 7="IDFC"; 8="IDFC"; 9="MUTHOOT"; 10="IDFC"; 11="IDFC"; 12="IDFC"; 13="IDFC"; 14="IDFC"; 15="IDFC" }
$department = @{ 2="TW"; 3="HL"; 4="TW"; 5="TW"; 6="TW"; 7="HL"; 8="HL"; 9="TW"; 10="HL"; 11="HL"; 12="HL"; 13="HL"; 14="HL"; 15="HL" }

foreach ($row in 2..15) {
    $ws.Cells.Item($row, 6).Value = $managerId[$row]
    $ws.Cells.Item($row, 10).Value = $process[$row]
    $ws.Cells.Item($row, 11).Value = $department[$row]
}
